$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1515941.8
$ws.Range("J17").Value = 1563299.2
$ws.Range("L17").Value = 4689897.6
$ws.Range("N17").Value = -4690233.6
$ws.Range("H125").Value = 693539.5
$ws.Range("I125").Value = 1703.3334
$ws.Range("J125").Value = 1471855.1
$ws.Range("K125").Value = 15330.0006
$ws.Range("L125").Value = 13246695.9
$ws.Range("M125").Value = -12870.0006
$ws.Range("N125").Value = -13251615.9
$ws.Range("H131").Value = 4529.469
$ws.Range("I131").Value = 275.375
$ws.Range("J131").Value = 4907.6113
$ws.Range("K131").Value = 826.125
$ws.Range("L131").Value = 14722.8339
$ws.Range("M131").Value = 4213.875
$ws.Range("N131").Value = -24802.8339
$ws.Range("H132").Value = 26510.863
$ws.Range("I132").Value = 29120.39
$ws.Range("J132").Value = 2503.2
$ws.Range("K132").Value = 87361.17
$ws.Range("L132").Value = 7509.599999999999
$ws.Range("M132").Value = -84831.17
$ws.Range("N132").Value = -12569.6
$ws.Range("H134").Value = 40056
$ws.Range("J134").Value = 40084.445
$ws.Range("L134").Value = 40084.445
$ws.Range("N134").Value = -50224.445
$ws.Range("H135").Value = 540.9778
$ws.Range("I135").Value = 540.9778
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4868.8002
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2333.8002
$ws.Range("H138").Value = 2826.37
$ws.Range("I138").Value = 1449.9108
$ws.Range("J138").Value = 4578.227
$ws.Range("K138").Value = 4349.732400000001
$ws.Range("L138").Value = 13734.681
$ws.Range("M138").Value = 790.2675999999992
$ws.Range("N138").Value = -24014.681
$ws.Range("N135").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2447
$ws.Range("I122").Value = 1716
$ws.Range("J122").Value = 3299.8333
$ws.Range("K122").Value = 5148
$ws.Range("L122").Value = 9899.499899999999
$ws.Range("M122").Value = -2698
$ws.Range("N122").Value = -14799.4999
$ws.Range("H132").Value = 1746.4615
$ws.Range("I132").Value = 906.9032
$ws.Range("J132").Value = 4999.75
$ws.Range("K132").Value = 2720.7096
$ws.Range("L132").Value = 14999.25
$ws.Range("M132").Value = -190.7096000000001
$ws.Range("N132").Value = -20059.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1383.1471
$ws.Range("I134").Value = 1100.4642
$ws.Range("J134").Value = 2702.3333
$ws.Range("K134").Value = 3301.3926
$ws.Range("L134").Value = 8106.999899999999
$ws.Range("M134").Value = -766.3925999999997
$ws.Range("N134").Value = -13176.9999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3271.2104
$ws.Range("I31").Value = 1784.8667
$ws.Range("J31").Value = 4240.5654
$ws.Range("K31").Value = 1784.8667
$ws.Range("L31").Value = 4240.5654
$ws.Range("M31").Value = -1489.8667
$ws.Range("N31").Value = -4830.5654
$ws.Range("H34").Value = 3271.2104
$ws.Range("I34").Value = 1784.8667
$ws.Range("J34").Value = 4240.5654
$ws.Range("K34").Value = 1784.8667
$ws.Range("L34").Value = 4240.5654
$ws.Range("M34").Value = -1582.8667
$ws.Range("N34").Value = -4644.5654
$ws.Range("H132").Value = 2163.6177
$ws.Range("I132").Value = 1535.7241
$ws.Range("J132").Value = 5805.4
$ws.Range("K132").Value = 4607.1723
$ws.Range("L132").Value = 17416.2
$ws.Range("M132").Value = -2077.1723
$ws.Range("N132").Value = -22476.2
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 1699.4906
$ws.Range("I134").Value = 1646.2
$ws.Range("J134").Value = 1863.4615
$ws.Range("K134").Value = 4938.6
$ws.Range("L134").Value = 5590.3845
$ws.Range("M134").Value = -2403.6
$ws.Range("N134").Value = -10660.3845
$ws.Range("H135").Value = 43949.715
$ws.Range("J135").Value = 43949.715
$ws.Range("L135").Value = 43949.715
$ws.Range("N135").Value = -54089.715
$ws.Range("N133").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 399401.84
$ws.Range("I5").Value = 425.55
$ws.Range("J5").Value = 731882.0600000001
$ws.Range("K5").Value = 1276.65
$ws.Range("L5").Value = 2195646.18
$ws.Range("M5").Value = -1164.65
$ws.Range("N5").Value = -2195870.18
$ws.Range("H131").Value = 2992.232
$ws.Range("I131").Value = 546.0769
$ws.Range("J131").Value = 3560.0894
$ws.Range("K131").Value = 1638.2307
$ws.Range("L131").Value = 10680.2682
$ws.Range("M131").Value = 3401.7693
$ws.Range("N131").Value = -20760.2682
$ws.Range("H133").Value = 2861.875
$ws.Range("I133").Value = 3532.2222
$ws.Range("K133").Value = 10596.6666
$ws.Range("M133").Value = -5536.6666
$ws.Range("H135").Value = 399401.84
$ws.Range("I135").Value = 425.55
$ws.Range("J135").Value = 731882.0600000001
$ws.Range("K135").Value = 3829.95
$ws.Range("L135").Value = 6586938.540000001
$ws.Range("M135").Value = -1294.95
$ws.Range("N135").Value = -6592008.540000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1858.8096
$ws.Range("I132").Value = 1510.7097
$ws.Range("J132").Value = 2839.818
$ws.Range("K132").Value = 4532.1291
$ws.Range("L132").Value = 8519.454000000002
$ws.Range("M132").Value = -2002.1291
$ws.Range("N132").Value = -13579.454
$ws.Range("H135").Value = 30838.889
$ws.Range("J135").Value = 30838.889
$ws.Range("L135").Value = 30838.889
$ws.Range("N135").Value = -40978.889

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4490.021
$ws.Range("I132").Value = 4216.8228
$ws.Range("J132").Value = 5759.5884
$ws.Range("K132").Value = 12650.4684
$ws.Range("L132").Value = 17278.7652
$ws.Range("M132").Value = -10120.4684
$ws.Range("N132").Value = -22338.7652
$ws.Range("H135").Value = 50172.867
$ws.Range("J135").Value = 50172.867
$ws.Range("L135").Value = 50172.867
$ws.Range("N135").Value = -60312.867
$ws.Range("H136").Value = 9010440
$ws.Range("I136").Value = 1470.75
$ws.Range("K136").Value = 4412.25
$ws.Range("M136").Value = -1862.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 30002.334
$ws.Range("J12").Value = 30002.334
$ws.Range("L12").Value = 30002.334
$ws.Range("N12").Value = -30286.334
$ws.Range("H132").Value = 1903.8478
$ws.Range("I132").Value = 1617.3077
$ws.Range("J132").Value = 3500.2856
$ws.Range("K132").Value = 4851.9231
$ws.Range("L132").Value = 10500.8568
$ws.Range("M132").Value = -2321.9231
$ws.Range("N132").Value = -15560.8568

Write-Host "Edit complete: updated Leve profit data across all sheets."
